$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1: drop left+right border from the thin box style -> borderId 4 (top+bottom only)
$c1_c1 = $ws1.Cells.Item(1, 3)
$c1_c1.Borders.Item(7).LineStyle = -4142
$c1_c1.Borders.Item(10).LineStyle = -4142

# D1: drop left border only from the thin box style -> borderId 5 (top+bottom+right)
$c1_d1 = $ws1.Cells.Item(1, 4)
$c1_d1.Borders.Item(7).LineStyle = -4142

# C2: rename header "fedcore" -> "approach"
$ws1.Cells.Item(2, 3).Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# C1: drop left+right border -> borderId 4
$c2_c1 = $ws2.Cells.Item(1, 3)
$c2_c1.Borders.Item(7).LineStyle = -4142
$c2_c1.Borders.Item(10).LineStyle = -4142

# D1: drop left border only -> borderId 5
$c2_d1 = $ws2.Cells.Item(1, 4)
$c2_d1.Borders.Item(7).LineStyle = -4142

# F1: drop left+right border -> borderId 4
$c2_f1 = $ws2.Cells.Item(1, 6)
$c2_f1.Borders.Item(7).LineStyle = -4142
$c2_f1.Borders.Item(10).LineStyle = -4142

# G1: drop left border only -> borderId 5
$c2_g1 = $ws2.Cells.Item(1, 7)
$c2_g1.Borders.Item(7).LineStyle = -4142

# C2 / F2: rename header "fedcore" -> "approach"
$ws2.Cells.Item(2, 3).Value = "approach"
$ws2.Cells.Item(2, 6).Value = "approach"

# G5: empty placeholder cell removed entirely
$ws2.Cells.Item(5, 7).ClearContents()
